$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the B16 cell formula to AVERAGE(B2:B14), matching the pattern of C16/D16/E16
$ws.Range("B16").Formula = "=AVERAGE(B2:B14)"

# Update the active selection cell on the sheet (as recorded when the file was last saved)
$ws.Range("B19").Select()
